# New weekly price record for "Pepino ensalada" at Vega Monumental Concepción.
# The new observation belongs at the top of this subset's date-ordered block
# (row 185), so insert a fresh row there and push the existing rows 185:196
# down to 186:197 (the sheet's dimension grows from A1:R196 to A1:R197).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(185).Insert()

$ws.Range("A185").Value = 11
$ws.Range("B185").Value = "Vega Monumental Concepción"
$ws.Range("C185").Value = "Bíobío"
$ws.Range("D185").Value = 45013
$ws.Range("E185").Value = 8
$ws.Range("F185").Value = 100112043
$ws.Range("G185").Value = "Pepino ensalada"
$ws.Range("H185").Value = "Sin especificar"
$ws.Range("I185").Value = "Primera"
$ws.Range("J185").Value = 250
$ws.Range("K185").Value = 6000
$ws.Range("L185").Value = 6500
$ws.Range("M185").Value = 6300
$ws.Range("N185").Value = "$/caja 60 unidades"
$ws.Range("O185").Value = "Región de Arica y Parinacota"
$ws.Range("P185").Value = 105
$ws.Range("Q185").Value = 60
$ws.Range("R185").Value = "Hortaliza"
